# Fixed string connector ID bug
# The socket_info sheet stored the connector/cluster id for some rows as a
# plain number (1, 3) instead of the expected string value. Update those
# cells to hold string values instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("socket_info")

# Row 2 (cu_id column A) was numeric 1 -> should be the string "abc"
$ws.Range("A2").Value = "abc"

# Row 4 (cu_id column A) was numeric 3 -> should be the string "c"
$ws.Range("A4").Value = "c"

# Update the active selection to mirror the saved workbook state
$ws.Range("I8:J8").Select()
